# Applies the "Add files via upload" edit to the tareas workbook:
#   - D3: fix the garbled/mojibake description text
#   - E3: Estado "En curso" -> "Completada"
#   - F3: % Avance 60 -> 100
#   - E15: Estado "En curso" -> "Pausada"
#   - G18: Prioridad "Media" -> "Alta"
#   - E22: Estado "En curso" -> "Pausada"
#   - F31: % Avance 10 -> 0
#   - Update the active selection / scroll position to D32 (topLeftCell A10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tareas")

$ws.Range("D3").Value = "Recuperacion de Trailer de 3 ejes, EDS Caldas Viejo"
$ws.Range("E3").Value = "Completada"
$ws.Range("F3").Value = 100

$ws.Range("E15").Value = "Pausada"

$ws.Range("G18").Value = "Alta"

$ws.Range("E22").Value = "Pausada"

$ws.Range("F31").Value = 0

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("D32").Select()
